$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to remain text,
# matching the inline-string storage used in the source file,
# so values like "0.998" or "608.19" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '69.755.77'
$ws.Cells.Item(2, 5).Value = '  +0.55%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.521.67'
$ws.Cells.Item(3, 5).Value = '  +0.96%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '608.19'
$ws.Cells.Item(5, 5).Value = '  -0.31%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '196.11'
$ws.Cells.Item(6, 5).Value = '  +5.40%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.631'
$ws.Cells.Item(7, 5).Value = '  +0.81%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -6.27%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '0.648'
$ws.Cells.Item(10, 5).Value = '  -0.30%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '53.72'
$ws.Cells.Item(11, 5).Value = '  +1.24%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.0000301'
$ws.Cells.Item(12, 5).Value = '  -2.01%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -0.33%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '4.080.57'
$ws.Cells.Item(14, 5).Value = '  +1.14%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '596.50'
$ws.Cells.Item(15, 5).Value = '  -0.71%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'Chainlink'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(16, 4).Value = '19.18'
$ws.Cells.Item(16, 5).Value = '  +1.74%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'Uniswap'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(17, 4).Value = '12.80'
$ws.Cells.Item(17, 5).Value = '  +1.66%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '69.904.02'
$ws.Cells.Item(18, 5).Value = '  +0.63%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '3.526.01'
$ws.Cells.Item(19, 5).Value = '  +0.94%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +1.48%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '0.991'
$ws.Cells.Item(21, 5).Value = '  +0.59%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '18.09'
$ws.Cells.Item(22, 5).Value = '  +6.01%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '5.34'
$ws.Cells.Item(23, 5).Value = '  +5.65%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '102.24'
$ws.Cells.Item(24, 5).Value = '  -2.75%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.55%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +4.32%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '10.88'
$ws.Cells.Item(27, 5).Value = '  -0.54%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '9.60'
$ws.Cells.Item(28, 5).Value = '  -1.01%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '33.53'
$ws.Cells.Item(29, 5).Value = '  -0.14%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +1.06%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +1.64%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '12.43'
$ws.Cells.Item(32, 5).Value = '  +0.16%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '0.115'
$ws.Cells.Item(33, 5).Value = '  +0.25%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '63.14'

# Row 35
$ws.Cells.Item(35, 4).Value = '0.0₃0847'
$ws.Cells.Item(35, 5).Value = '  +8.95%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '3.713.21'
$ws.Cells.Item(36, 5).Value = '  +2.82%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Fetch.AI'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(37, 4).Value = '3.08'
$ws.Cells.Item(37, 5).Value = '  -3.44%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Dai'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(38, 4).Value = '0.999'
$ws.Cells.Item(38, 5).Value = '  +0.15%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.06%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -1.19%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '36.56'
$ws.Cells.Item(41, 5).Value = '  -0.39%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '486.05'
$ws.Cells.Item(42, 5).Value = '  -7.31%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -3.73%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -0.62%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -1.28%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '2.82'
$ws.Cells.Item(46, 5).Value = '  -3.82%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -1.29%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.26%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '8.45'
$ws.Cells.Item(49, 5).Value = '  -3.90%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +1.49%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +10.88%  '
